$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.8740398169273846, 1.2329050188459543, -0.030582831213484627, 2.5718628340153167)"
$ws.Range("C2").Value = "NIG(1.2854475829052339, 0.9380794127144989, 3.9209204097618917, 6.16521099304507)"
$ws.Range("D2").Value = "NIG(0.8140537683887275, 0.5193930603136606, 1.7078718547586869, 3.038332692316259)"
$ws.Range("E2").Value = "JSU(-1.4974458139227005, 1.3217501917493517, 1.9127348890913867, 5.178704868676462)"
